$d = $word.ActiveDocument

$d.Content.Find.Execute("ĐỒ ÁN", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ĐỒ ÁN LẦN 1", 2)
